# This script updates the practice-sheet date heading and all 100
# arithmetic answers in the single table, cell by cell, matching the
# target OOXML diff exactly. We assign Range.Text directly (rather than
# Find.Execute) because several cells share identical old text (e.g.
# "58+29=87" appears twice) and Find.Execute in this runtime searches
# from the start of the story instead of respecting the scoped Range,
# which would otherwise edit the wrong (first) occurrence.
$d = $word.ActiveDocument

# Update the date heading paragraph
$dateRange = $d.Paragraphs.Item(1).Range
$dateRange.Text = "2023-11-05 Sunday"

# Update the table cells (20 rows x 5 columns)
$t = $d.Tables.Item(1)

$cell = $t.Cell(1, 1)
$cell.Range.Text = "9+23=32"  # was: 70-57=13
$cell = $t.Cell(1, 2)
$cell.Range.Text = "91-57=34"  # was: 93-57=36
$cell = $t.Cell(1, 3)
$cell.Range.Text = "64-16=48"  # was: 19+23=42
$cell = $t.Cell(1, 4)
$cell.Range.Text = "36+15=51"  # was: 81-56=25
$cell = $t.Cell(1, 5)
$cell.Range.Text = "45+49=94"  # was: 23+9=32

$cell = $t.Cell(2, 1)
$cell.Range.Text = "40-31=9"  # was: 46+47=93
$cell = $t.Cell(2, 2)
$cell.Range.Text = "70-4=66"  # was: 7+85=92
$cell = $t.Cell(2, 3)
$cell.Range.Text = "90-84=6"  # was: 58+29=87
$cell = $t.Cell(2, 4)
$cell.Range.Text = "47+36=83"  # was: 27+39=66
$cell = $t.Cell(2, 5)
$cell.Range.Text = "40-1=39"  # was: 90-3=87

$cell = $t.Cell(3, 1)
$cell.Range.Text = "27+68=95"  # was: 48+16=64
$cell = $t.Cell(3, 2)
$cell.Range.Text = "4+79=83"  # was: 72-26=46
$cell = $t.Cell(3, 3)
$cell.Range.Text = "53-24=29"  # was: 93-35=58
$cell = $t.Cell(3, 4)
$cell.Range.Text = "91-78=13"  # was: 53-35=18
$cell = $t.Cell(3, 5)
$cell.Range.Text = "7+77=84"  # was: 64-45=19

$cell = $t.Cell(4, 1)
$cell.Range.Text = "72-66=6"  # was: 60-26=34
$cell = $t.Cell(4, 2)
$cell.Range.Text = "4+48=52"  # was: 4+87=91
$cell = $t.Cell(4, 3)
$cell.Range.Text = "74-55=19"  # was: 37-8=29
$cell = $t.Cell(4, 4)
$cell.Range.Text = "48+14=62"  # was: 27+8=35
$cell = $t.Cell(4, 5)
$cell.Range.Text = "48+28=76"  # was: 30-4=26

$cell = $t.Cell(5, 1)
$cell.Range.Text = "93-9=84"  # was: 19+68=87
$cell = $t.Cell(5, 2)
$cell.Range.Text = "37+47=84"  # was: 27+35=62
$cell = $t.Cell(5, 3)
$cell.Range.Text = "32-5=27"  # was: 72-19=53
$cell = $t.Cell(5, 4)
$cell.Range.Text = "33+48=81"  # was: 13+69=82
$cell = $t.Cell(5, 5)
$cell.Range.Text = "92-85=7"  # was: 72-37=35

$cell = $t.Cell(6, 1)
$cell.Range.Text = "12+39=51"  # was: 50-43=7
$cell = $t.Cell(6, 2)
$cell.Range.Text = "9+3=12"  # was: 38+48=86
$cell = $t.Cell(6, 3)
$cell.Range.Text = "63-56=7"  # was: 19+8=27
$cell = $t.Cell(6, 4)
$cell.Range.Text = "63-35=28"  # was: 17+68=85
$cell = $t.Cell(6, 5)
$cell.Range.Text = "43-27=16"  # was: 71-35=36

$cell = $t.Cell(7, 1)
$cell.Range.Text = "30-4=26"  # was: 73-8=65
$cell = $t.Cell(7, 2)
$cell.Range.Text = "73-34=39"  # was: 9+4=13
$cell = $t.Cell(7, 3)
$cell.Range.Text = "35-26=9"  # was: 28+14=42
$cell = $t.Cell(7, 4)
$cell.Range.Text = "56+27=83"  # was: 64+8=72
$cell = $t.Cell(7, 5)
$cell.Range.Text = "75-56=19"  # was: 7+47=54

$cell = $t.Cell(8, 1)
$cell.Range.Text = "71-26=45"  # was: 73-44=29
$cell = $t.Cell(8, 2)
$cell.Range.Text = "63-39=24"  # was: 33+28=61
$cell = $t.Cell(8, 3)
$cell.Range.Text = "46-29=17"  # was: 64-27=37
$cell = $t.Cell(8, 4)
$cell.Range.Text = "47-19=28"  # was: 64-9=55
$cell = $t.Cell(8, 5)
$cell.Range.Text = "74+7=81"  # was: 55+38=93

$cell = $t.Cell(9, 1)
$cell.Range.Text = "73-27=46"  # was: 74-28=46
$cell = $t.Cell(9, 2)
$cell.Range.Text = "93-45=48"  # was: 6+6=12
$cell = $t.Cell(9, 3)
$cell.Range.Text = "39+42=81"  # was: 91-35=56
$cell = $t.Cell(9, 4)
$cell.Range.Text = "48+33=81"  # was: 16+67=83
$cell = $t.Cell(9, 5)
$cell.Range.Text = "50-45=5"  # was: 83-5=78

$cell = $t.Cell(10, 1)
$cell.Range.Text = "8+47=55"  # was: 92-86=6
$cell = $t.Cell(10, 2)
$cell.Range.Text = "49+29=78"  # was: 35+6=41
$cell = $t.Cell(10, 3)
$cell.Range.Text = "59+39=98"  # was: 43+38=81
$cell = $t.Cell(10, 4)
$cell.Range.Text = "9+46=55"  # was: 8+75=83
$cell = $t.Cell(10, 5)
$cell.Range.Text = "66-49=17"  # was: 9+7=16

$cell = $t.Cell(11, 1)
$cell.Range.Text = "26+45=71"  # was: 86-18=68
$cell = $t.Cell(11, 2)
$cell.Range.Text = "39+58=97"  # was: 58+29=87
$cell = $t.Cell(11, 3)
$cell.Range.Text = "71-39=32"  # was: 84-58=26
$cell = $t.Cell(11, 4)
$cell.Range.Text = "25+37=62"  # was: 57+4=61
$cell = $t.Cell(11, 5)
$cell.Range.Text = "45+39=84"  # was: 16+17=33

$cell = $t.Cell(12, 1)
$cell.Range.Text = "55+37=92"  # was: 92-8=84
$cell = $t.Cell(12, 2)
$cell.Range.Text = "9+38=47"  # was: 48+48=96
$cell = $t.Cell(12, 3)
$cell.Range.Text = "34+57=91"  # was: 4+67=71
$cell = $t.Cell(12, 4)
$cell.Range.Text = "41-13=28"  # was: 13+38=51
$cell = $t.Cell(12, 5)
$cell.Range.Text = "8+6=14"  # was: 80-26=54

$cell = $t.Cell(13, 1)
$cell.Range.Text = "70-22=48"  # was: 54-45=9
$cell = $t.Cell(13, 2)
$cell.Range.Text = "60-15=45"  # was: 19+67=86
$cell = $t.Cell(13, 3)
$cell.Range.Text = "28+4=32"  # was: 96-29=67
$cell = $t.Cell(13, 4)
$cell.Range.Text = "19+32=51"  # was: 42-17=25
$cell = $t.Cell(13, 5)
$cell.Range.Text = "33-16=17"  # was: 46-19=27

$cell = $t.Cell(14, 1)
$cell.Range.Text = "26+58=84"  # was: 57+19=76
$cell = $t.Cell(14, 2)
$cell.Range.Text = "94-25=69"  # was: 19+76=95
$cell = $t.Cell(14, 3)
$cell.Range.Text = "91-29=62"  # was: 9+14=23
$cell = $t.Cell(14, 4)
$cell.Range.Text = "28+26=54"  # was: 17+68=85
$cell = $t.Cell(14, 5)
$cell.Range.Text = "40-34=6"  # was: 47+5=52

$cell = $t.Cell(15, 1)
$cell.Range.Text = "42+19=61"  # was: 19+6=25
$cell = $t.Cell(15, 2)
$cell.Range.Text = "58+13=71"  # was: 45-39=6
$cell = $t.Cell(15, 3)
$cell.Range.Text = "29+58=87"  # was: 78+17=95
$cell = $t.Cell(15, 4)
$cell.Range.Text = "30-17=13"  # was: 39+17=56
$cell = $t.Cell(15, 5)
$cell.Range.Text = "33-18=15"  # was: 81-46=35

$cell = $t.Cell(16, 1)
$cell.Range.Text = "62-44=18"  # was: 5+87=92
$cell = $t.Cell(16, 2)
$cell.Range.Text = "15+47=62"  # was: 64-28=36
$cell = $t.Cell(16, 3)
$cell.Range.Text = "28+47=75"  # was: 17+46=63
$cell = $t.Cell(16, 4)
$cell.Range.Text = "55+18=73"  # was: 4+67=71
$cell = $t.Cell(16, 5)
$cell.Range.Text = "40-13=27"  # was: 29+64=93

$cell = $t.Cell(17, 1)
$cell.Range.Text = "40-34=6"  # was: 53-44=9
$cell = $t.Cell(17, 2)
$cell.Range.Text = "85-38=47"  # was: 71-34=37
$cell = $t.Cell(17, 3)
$cell.Range.Text = "41-15=26"  # was: 8+89=97
$cell = $t.Cell(17, 4)
$cell.Range.Text = "25+9=34"  # was: 16+48=64
$cell = $t.Cell(17, 5)
$cell.Range.Text = "67+14=81"  # was: 32-17=15

$cell = $t.Cell(18, 1)
$cell.Range.Text = "44+48=92"  # was: 63-26=37
$cell = $t.Cell(18, 2)
$cell.Range.Text = "72-4=68"  # was: 12+9=21
$cell = $t.Cell(18, 3)
$cell.Range.Text = "80-75=5"  # was: 64-39=25
$cell = $t.Cell(18, 4)
$cell.Range.Text = "46-28=18"  # was: 42+9=51
$cell = $t.Cell(18, 5)
$cell.Range.Text = "57-28=29"  # was: 84-79=5

$cell = $t.Cell(19, 1)
$cell.Range.Text = "47-9=38"  # was: 26+46=72
$cell = $t.Cell(19, 2)
$cell.Range.Text = "63+28=91"  # was: 47+29=76
$cell = $t.Cell(19, 3)
$cell.Range.Text = "18+46=64"  # was: 68+3=71
$cell = $t.Cell(19, 4)
$cell.Range.Text = "56+37=93"  # was: 71-15=56
$cell = $t.Cell(19, 5)
$cell.Range.Text = "81-24=57"  # was: 34-6=28

$cell = $t.Cell(20, 1)
$cell.Range.Text = "38+36=74"  # was: 88+3=91
$cell = $t.Cell(20, 2)
$cell.Range.Text = "25-16=9"  # was: 26-17=9
$cell = $t.Cell(20, 3)
$cell.Range.Text = "48+28=76"  # was: 78-49=29
$cell = $t.Cell(20, 4)
$cell.Range.Text = "72-28=44"  # was: 87+4=91
$cell = $t.Cell(20, 5)
$cell.Range.Text = "26+18=44"  # was: 84-16=68

Write-Host "Done updating document."
